$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A56").Value = 46005
$ws.Range("B56").Value = 121
$ws.Range("C56").Value = 135
$ws.Range("D56").Value = 126

$ws.Range("A56").NumberFormat = $ws.Range("A55").NumberFormat
